$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the grade values in column C
$ws.Range("C18").Value = 86
$ws.Range("C19").Value = 88
$ws.Range("C20").Value = 89
$ws.Range("C21").Value = 86
$ws.Range("C22").Value = 85

# Update the active cell selection to C22
$ws.Activate()
$ws.Range("C22").Select()
